$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.852944
$ws.Range("H2").Value = 83.558832
$ws.Range("I2").Value = 0.2559209115167818
$ws.Range("J2").Value = 0.2559209115167818
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 72.36658138501866
$ws.Range("R2").Value = 651.299232465168
$ws.Range("S2").Value = 0.08870252885810061
$ws.Range("T2").Value = 0.08870252885810061
$ws.Range("G3").Value = 27.852944
$ws.Range("H3").Value = 83.558832
$ws.Range("I3").Value = 0.2559209115167818
$ws.Range("J3").Value = 0.2559209115167818
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 120.6980403727466
$ws.Range("R3").Value = 1086.28236335472
$ws.Range("S3").Value = 0.147944274890069
$ws.Range("T3").Value = 0.147944274890069
$ws.Range("G4").Value = 27.852944
$ws.Range("H4").Value = 83.558832
$ws.Range("I4").Value = 0.2559209115167818
$ws.Range("J4").Value = 0.2559209115167818
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 13.06872202089067
$ws.Range("R4").Value = 117.618498188016
$ws.Range("S4").Value = 0.01601884005034032
$ws.Range("T4").Value = 0.01601884005034032
$ws.Range("G5").Value = 27.852944
$ws.Range("H5").Value = 83.558832
$ws.Range("I5").Value = 0.2559209115167818
$ws.Range("J5").Value = 0.2559209115167818
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 2.655759641770666
$ws.Range("R5").Value = 23.901836775936
$ws.Range("S5").Value = 0.003255267718271818
$ws.Range("T5").Value = 0.003255267718271817
$ws.Range("I6").Value = 0.3112048767201538
$ws.Range("J6").Value = 0.3112048767201538
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 87.99919047297911
$ws.Range("R6").Value = 791.9927142568121
$ws.Range("S6").Value = 0.1078640248444135
$ws.Range("T6").Value = 0.1078640248444135
$ws.Range("I7").Value = 0.3112048767201538
$ws.Range("J7").Value = 0.3112048767201538
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("S7").Value = 0.1799031566265634
$ws.Range("T7").Value = 0.1799031566265634
$ws.Range("I8").Value = 0.3112048767201538
$ws.Range("J8").Value = 0.3112048767201538
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 15.89182377202711
$ws.Range("R8").Value = 143.026413948244
$ws.Range("S8").Value = 0.01947922549009492
$ws.Range("T8").Value = 0.01947922549009492
$ws.Range("I9").Value = 0.3112048767201538
$ws.Range("J9").Value = 0.3112048767201538
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 3.229456112113777
$ws.Range("R9").Value = 29.065105009024
$ws.Range("S9").Value = 0.003958469759081987
$ws.Range("T9").Value = 0.003958469759081986
$ws.Range("G10").Value = 30.14135433333333
$ws.Range("H10").Value = 90.42406299999999
$ws.Range("I10").Value = 0.2769474880406526
$ws.Range("J10").Value = 0.2769474880406526
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 78.31225206993743
$ws.Range("R10").Value = 704.8102686294369
$ws.Range("S10").Value = 0.0959903682919384
$ws.Range("T10").Value = 0.09599036829193838
$ws.Range("G11").Value = 30.14135433333333
$ws.Range("H11").Value = 90.42406299999999
$ws.Range("I11").Value = 0.2769474880406526
$ws.Range("J11").Value = 0.2769474880406526
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 130.6146453392477
$ws.Range("R11").Value = 1175.53180805323
$ws.Range("S11").Value = 0.1600994426675198
$ws.Range("T11").Value = 0.1600994426675198
$ws.Range("G12").Value = 30.14135433333333
$ws.Range("H12").Value = 90.42406299999999
$ws.Range("I12").Value = 0.2769474880406526
$ws.Range("J12").Value = 0.2769474880406526
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 14.14245406573544
$ws.Range("R12").Value = 127.282086591619
$ws.Range("S12").Value = 0.01733495511161401
$ws.Range("T12").Value = 0.01733495511161401
$ws.Range("G13").Value = 30.14135433333333
$ws.Range("H13").Value = 90.42406299999999
$ws.Range("I13").Value = 0.2769474880406526
$ws.Range("J13").Value = 0.2769474880406526
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 2.873958041447111
$ws.Range("R13").Value = 25.865622373024
$ws.Range("S13").Value = 0.003522721969580392
$ws.Range("T13").Value = 0.003522721969580391
$ws.Range("G14").Value = 16.970158
$ws.Range("H14").Value = 50.910474
$ws.Range("I14").Value = 0.1559267237224118
$ws.Range("J14").Value = 0.1559267237224118
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 44.09129318694734
$ws.Range("R14").Value = 396.821638682526
$ws.Range("S14").Value = 0.05404441016079044
$ws.Range("T14").Value = 0.05404441016079043
$ws.Range("G15").Value = 16.970158
$ws.Range("H15").Value = 50.910474
$ws.Range("I15").Value = 0.1559267237224118
$ws.Range("J15").Value = 0.1559267237224118
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 73.53853924439333
$ws.Range("R15").Value = 661.84685319954
$ws.Range("S15").Value = 0.09013904311443358
$ws.Range("T15").Value = 0.09013904311443356
$ws.Range("G16").Value = 16.970158
$ws.Range("H16").Value = 50.910474
$ws.Range("I16").Value = 0.1559267237224118
$ws.Range("J16").Value = 0.1559267237224118
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 7.962471670951334
$ws.Range("R16").Value = 71.66224503856201
$ws.Range("S16").Value = 0.009759910716475905
$ws.Range("T16").Value = 0.009759910716475903
$ws.Range("G17").Value = 16.970158
$ws.Range("H17").Value = 50.910474
$ws.Range("I17").Value = 0.1559267237224118
$ws.Range("J17").Value = 0.1559267237224118
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 1.618093251861333
$ws.Range("R17").Value = 14.562839266752
$ws.Range("S17").Value = 0.00198335973071185
$ws.Range("T17").Value = 0.00198335973071185
